# Update the noise trajectory values in columns B and C (rows 2-9)
# on the active worksheet, matching the newly generated random trajectory.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "[0, 0, 0]"
$ws.Range("C2").Value = "[0, 1, 0]"

$ws.Range("B3").Value = "[0, 1, 0]"
$ws.Range("C3").Value = "[1, 1, 1]"

$ws.Range("B4").Value = "[1, 1, 1]"
$ws.Range("C4").Value = "[0, 0, 0]"

$ws.Range("B5").Value = "[1, 0, 0]"
$ws.Range("C5").Value = "[0, 1, 0]"

$ws.Range("B6").Value = "[1, 1, 0]"
$ws.Range("C6").Value = "[1, 1, 0]"

$ws.Range("B7").Value = "[0, 0, 1]"
$ws.Range("C7").Value = "[1, 1, 0]"

$ws.Range("B8").Value = "[1, 0, 1]"
$ws.Range("C8").Value = "[1, 1, 0]"

$ws.Range("B9").Value = "[0, 1, 1]"
$ws.Range("C9").Value = "[0, 0, 1]"
